$d = $word.ActiveDocument

$d.Content.Find.Execute("901÷7=128, 5", $true, $false, $false, $false, $false, $true, 1, $false, "496÷2=248, 0", 2) | Out-Null
$d.Content.Find.Execute("717÷2=358, 1", $true, $false, $false, $false, $false, $true, 1, $false, "253÷4=63, 1", 2) | Out-Null
$d.Content.Find.Execute("892÷7=127, 3", $true, $false, $false, $false, $false, $true, 1, $false, "283÷6=47, 1", 2) | Out-Null
$d.Content.Find.Execute("158÷3=52, 2", $true, $false, $false, $false, $false, $true, 1, $false, "244÷7=34, 6", 2) | Out-Null
$d.Content.Find.Execute("772÷2=386, 0", $true, $false, $false, $false, $false, $true, 1, $false, "993÷5=198, 3", 2) | Out-Null
$d.Content.Find.Execute("457÷9=50, 7", $true, $false, $false, $false, $false, $true, 1, $false, "577÷3=192, 1", 2) | Out-Null
$d.Content.Find.Execute("419÷3=139, 2", $true, $false, $false, $false, $false, $true, 1, $false, "753÷8=94, 1", 2) | Out-Null
$d.Content.Find.Execute("602÷7=86, 0", $true, $false, $false, $false, $false, $true, 1, $false, "337÷9=37, 4", 2) | Out-Null
$d.Content.Find.Execute("512÷3=170, 2", $true, $false, $false, $false, $false, $true, 1, $false, "440÷8=55, 0", 2) | Out-Null
$d.Content.Find.Execute("906÷7=129, 3", $true, $false, $false, $false, $false, $true, 1, $false, "741÷5=148, 1", 2) | Out-Null
$d.Content.Find.Execute("123÷2=61, 1", $true, $false, $false, $false, $false, $true, 1, $false, "646÷6=107, 4", 2) | Out-Null
$d.Content.Find.Execute("978÷4=244, 2", $true, $false, $false, $false, $false, $true, 1, $false, "131÷2=65, 1", 2) | Out-Null
$d.Content.Find.Execute("615÷9=68, 3", $true, $false, $false, $false, $false, $true, 1, $false, "371÷6=61, 5", 2) | Out-Null
$d.Content.Find.Execute("214÷3=71, 1", $true, $false, $false, $false, $false, $true, 1, $false, "581÷4=145, 1", 2) | Out-Null
$d.Content.Find.Execute("173÷2=86, 1", $true, $false, $false, $false, $false, $true, 1, $false, "340÷8=42, 4", 2) | Out-Null
$d.Content.Find.Execute("973÷5=194, 3", $true, $false, $false, $false, $false, $true, 1, $false, "148÷5=29, 3", 2) | Out-Null
$d.Content.Find.Execute("313÷4=78, 1", $true, $false, $false, $false, $false, $true, 1, $false, "534÷6=89, 0", 2) | Out-Null
$d.Content.Find.Execute("284÷7=40, 4", $true, $false, $false, $false, $false, $true, 1, $false, "225÷7=32, 1", 2) | Out-Null
$d.Content.Find.Execute("606÷2=303, 0", $true, $false, $false, $false, $false, $true, 1, $false, "439÷8=54, 7", 2) | Out-Null
$d.Content.Find.Execute("847÷8=105, 7", $true, $false, $false, $false, $false, $true, 1, $false, "737÷7=105, 2", 2) | Out-Null
$d.Content.Find.Execute("908÷7=129, 5", $true, $false, $false, $false, $false, $true, 1, $false, "993÷3=331, 0", 2) | Out-Null
$d.Content.Find.Execute("501÷6=83, 3", $true, $false, $false, $false, $false, $true, 1, $false, "621÷8=77, 5", 2) | Out-Null
$d.Content.Find.Execute("278÷6=46, 2", $true, $false, $false, $false, $false, $true, 1, $false, "703÷3=234, 1", 2) | Out-Null
$d.Content.Find.Execute("198÷6=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "566÷4=141, 2", 2) | Out-Null
$d.Content.Find.Execute("615÷8=76, 7", $true, $false, $false, $false, $false, $true, 1, $false, "234÷2=117, 0", 2) | Out-Null
